$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# Field type for CreateDate (row 13) and LastUpdate (row 15) changed from DATE to TIMESTAMP
$ws.Range("D13").Value = "TIMESTAMP"
$ws.Range("D15").Value = "TIMESTAMP"

# Update the active selection to match the edited cell
$ws.Range("D15").Select()
